# UNTREC now get the QUOTED AMOUNT instead of INVOICING AMOUNT
#
# - loginQA: a second login (grover / Mal@0000) is added in row 3, mirroring
#   the existing row 2 formatting, plus a mailto hyperlink on the password cell.
# - PO_Detail: the quote/DA?/unit lookup values change, and a new
#   "In Service Date" column is appended.
# - Unit_to_Reconcile_Output: now resolves to a single reconciled unit (the
#   quoted-amount record) instead of two rows that included the posted /
#   invoiced PO amount.

function Set-TextValue($range, $value) {
    # Force the cell to be stored as text (shared string) even when the
    # value looks numeric/date-like (leading zeros, "$", "/", etc.) so we
    # don't silently turn it into a real number/date/currency cell.
    $range.NumberFormat = "@"
    $range.Value2 = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: loginQA - add a second login row (grover / Mal@0000)
# ---------------------------------------------------------------------
$wsLogin = $wb.Worksheets.Item("loginQA")

$wsLogin.Range("A2").Copy()
$wsLogin.Range("A3").PasteSpecial(-4122)   # xlPasteFormats
$wsLogin.Range("A3").Value2 = "grover"

$wsLogin.Range("B3").Value2 = "Mal@0000"
$wsLogin.Hyperlinks.Add($wsLogin.Range("B3"), "mailto:Mal@0000")

# Re-apply B2's formatting after the hyperlink is created - Hyperlinks.Add
# stamps its own (duplicate) style, so copy the real one over it last.
$wsLogin.Range("B2").Copy()
$wsLogin.Range("B3").PasteSpecial(-4122)   # xlPasteFormats

$wsLogin.Activate()
$wsLogin.Range("C17").Select()

# ---------------------------------------------------------------------
# Sheet 2: PO_Detail - new Quote / DA? / Unit lookup + In Service Date
# ---------------------------------------------------------------------
$wsPO = $wb.Worksheets.Item("PO_Detail")

$wsPO.Range("D1").Value2 = "In Service Date"

$wsPO.Range("A2").Value2 = "'369185"
$wsPO.Range("B2").Value2 = "'1"
$wsPO.Range("C2").Value2 = "'00995845"
$wsPO.Range("D2").Value2 = "'05/01/2017"

$wsPO.Activate()
$wsPO.Range("A2").Select()

# ---------------------------------------------------------------------
# Sheet 3: Unit_to_Reconcile_Output - rebuilt with a single reconciled
# row (quoted amount only, no posted/invoiced amount column populated).
# Deleting + re-adding also restores the sheet's internal id sequencing.
# ---------------------------------------------------------------------
$sheetName = $wb.Worksheets.Item("Unit_to_Reconcile_Output").Name
$afterSheet = $wb.Worksheets.Item("PO_Detail")
$wb.Worksheets.Item("Unit_to_Reconcile_Output").Delete()

$wsOut = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$wsOut.Name = $sheetName

$wsOut.Range("A1").Value2 = "Unit"
$wsOut.Range("B1").Value2 = "Unit_desc"
$wsOut.Range("C1").Value2 = "VIN"
$wsOut.Range("D1").Value2 = "Vendor"
$wsOut.Range("E1").Value2 = "PO Amount_Shown"
$wsOut.Range("F1").Value2 = "PON"
$wsOut.Range("G1").Value2 = "INV_Entered"
$wsOut.Range("H1").Value2 = "PO Amount_posted"

Set-TextValue $wsOut.Range("A2") "00998503"
$wsOut.Range("B2").Value2 = "00998503 - 2017 Chevrolet Express 2500 Work Van Rear-wheel Drive Cargo Van (CG23405)"
$wsOut.Range("C2").Value2 = "AAAAAAAAAAAAAAAAA"
$wsOut.Range("D2").Value2 = "00158470 - Al Piemonte Chevrolet"
Set-TextValue $wsOut.Range("E2") "$24,914.50"
$wsOut.Range("F2").Value2 = "PON00213850/0"
$wsOut.Range("G2").Value2 = "INV00213850"

# Re-activate PO_Detail so it remains the workbook's active tab.
$wsPO.Activate()
